# ---------------------------------------------------------------------------
# 2021AP.xlsx — "Adding work on datasets and models"
#
# Week 5 (rows 102-126) previously only had columns A-F filled in; this adds
# the remaining game-result columns G-O for each of those 25 teams.
#
# Week 6 results are then appended as 25 brand-new rows (127-151), pushing
# the trailing header/legend row down from row 127 to row 152.
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Make room for the 25 new Week 6 rows by inserting before the legend row
#    (old row 127). This shifts the legend row down to row 152 automatically,
#    exactly like selecting rows 127:151 in Excel and choosing Insert.
# ---------------------------------------------------------------------------
$ws.Rows("127:151").Insert()

# ---------------------------------------------------------------------------
# 2) Fill in the missing Week 5 columns (G:O) for rows 102-126.
#    G=Opp. Rank, H=Opp. P5, I=Home, J=Result (W/L), K=Points Scored,
#    L=Points Against, M=Margin(=K-L), N=Next Week Rank, O=Movement(=C-N)
# ---------------------------------------------------------------------------
$week5 = @(
    @(102, 12, 1, 1, "W", 42, 21, 1),
    @(103,  8, 1, 1, "W", 37,  0, 2),
    @(104, 26, 1, 0, "L", 24, 31, 8),
    @(105, 26, 1, 1, "W", 24,  0, 4),
    @(106, 26, 1, 0, "W", 51, 14, 3),
    @(107, 26, 1, 0, "W", 37, 31, 6),
    @(108,  9, 1, 0, "W", 24, 13, 5),
    @(109,  2, 1, 0, "L",  0, 37, 13),
    @(110,  7, 0, 1, "L", 13, 24, 14),
    @(111, 26, 1, 0, "L", 13, 20, 20),
    @(112, 26, 1, 0, "W", 52, 13, 7),
    @(113,  1, 1, 0, "L", 21, 42, 17),
    @(114, 26, 0, 0, "W", 34, 20, 10),
    @(115, 26, 1, 0, "W", 38, 17, 9),
    @(116, 26, 1, 1, "L", 22, 26, 26),
    @(117, 26, 0, 1, "W", 59,  6, 15),
    @(118, 26, 0, 1, "W", 48, 31, 11),
    @(119, 26, 0, 0, "L", 24, 27, 26),
    @(120, 21, 1, 1, "W", 17, 14, 12),
    @(121, 26, 1, 1, "L", 23, 42, 26),
    @(122, 19, 1, 0, "L", 14, 17, 26),
    @(123, 26, 1, 0, "W", 24, 19, 18),
    @(124, 26, 0, 1, "W", 34, 27, 23),
    @(125, 26, 1, 1, "W", 37, 34, 19),
    @(126, 26, 1, 1, "W", 19, 13, 26)
)

foreach ($row in $week5) {
    $r = $row[0]
    $ws.Cells.Item($r, 7).Value  = $row[1]   # G - Opp. Rank
    $ws.Cells.Item($r, 8).Value  = $row[2]   # H - Opp. P5
    $ws.Cells.Item($r, 9).Value  = $row[3]   # I - Home
    $ws.Cells.Item($r, 10).Value = $row[4]   # J - Result
    $ws.Cells.Item($r, 11).Value = $row[5]   # K - Points Scored
    $ws.Cells.Item($r, 12).Value = $row[6]   # L - Points Against
    $ws.Cells.Item($r, 13).Formula = "=K$r-L$r"   # M - Margin
    $ws.Cells.Item($r, 14).Value = $row[7]   # N - Next Week Rank
    $ws.Cells.Item($r, 15).Formula = "=C$r-N$r"   # O - Movement
}

# ---------------------------------------------------------------------------
# 3) Write the 25 new Week 6 rows (127-151), columns A-F only:
#    A=Team, B=Week, C=Rank (running count), D=W, E=L, F=Winning % (=D/(D+E))
# ---------------------------------------------------------------------------
$week6 = @(
    @(127, "Alabama",          5, 0),
    @(128, "Georgia",          5, 0),
    @(129, "Iowa",              5, 0),
    @(130, "Penn State",        5, 0),
    @(131, "Cincinnati",        4, 0),
    @(132, "Oklahoma",          5, 0),
    @(133, "Ohio State",        4, 1),
    @(134, "Oregon",            4, 1),
    @(135, "Michigan",          5, 0),
    @(136, "BYU",               5, 0),
    @(137, "Michigan State",    5, 0),
    @(138, "Oklahoma State",    5, 0),
    @(139, "Arkansas",          4, 1),
    @(140, "Notre Dame",        4, 1),
    @(141, "Coastal Carolina",  5, 0),
    @(142, "Kentucky",          5, 0),
    @(143, "Ole Miss",          3, 1),
    @(144, "Auburn",            4, 4),
    @(145, "Wake Forest",       5, 0),
    @(146, "Florida",           3, 2),
    @(147, "Texas",             4, 1),
    @(148, "Arizona State",     4, 1),
    @(149, "NC State",          4, 1),
    @(150, "SMU",               5, 0),
    @(151, "San Diego State",   4, 0)
)

foreach ($row in $week6) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]     # A - Team
    $ws.Cells.Item($r, 2).Value = 6           # B - Week
    if ($r -eq 127) {
        $ws.Cells.Item($r, 3).Value = 1       # C - Rank (first row is literal 1)
    } else {
        $prev = $r - 1
        $ws.Cells.Item($r, 3).Formula = "=C$prev+1"
    }
    $ws.Cells.Item($r, 4).Value = $row[2]     # D - W
    $ws.Cells.Item($r, 5).Value = $row[3]     # E - L
    $ws.Cells.Item($r, 6).Formula = "=D$r/(D$r+E$r)"   # F - Winning %
}

# ---------------------------------------------------------------------------
# 4) Update the saved view state: scrolled down with Q121 selected, matching
#    the author's position at the time of saving.
# ---------------------------------------------------------------------------
$ws.Application.GoTo($ws.Range("Q121"), $true)
$ws.Range("Q121").Select()
